$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'246.36"
$ws.Range("D3").Value = "'22.68"
$ws.Range("D4").Value = "'5.536"
$ws.Range("D5").Value = "'0.05605"
$ws.Range("B6").Value = 'KuCoinToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D6").Value = "'6.470"
$ws.Range("E6").Value = '5KuCoinTokenKCS'
$ws.Range("B7").Value = 'MXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D7").Value = "'0.8018"
$ws.Range("E7").Value = '6MXTokenMX'
$ws.Range("B8").Value = 'FTXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D8").Value = "'1.053"
$ws.Range("E8").Value = '7FTXTokenFTT'
$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9").Value = "'0.1422"
$ws.Range("E9").Value = '8WazirXWRX'
$ws.Range("B10").Value = 'MandalaExchangeToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D10").Value = "'0.07265"
$ws.Range("E10").Value = '9MandalaExchangeTokenMDX'
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").Value = "'0.03204"
$ws.Range("E11").Value = '10LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = "'0.02968"
$ws.Range("E12").Value = '11BitrueCoinBTR'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = "'0.09270"
$ws.Range("E13").Value = '12BitMartTokenBMX'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = "'0.001676"
$ws.Range("E14").Value = '13BitForexTokenBF'
$ws.Range("B15").Value = 'MCDex'
$ws.Range("C15").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D15").Value = "'2.975"
$ws.Range("E15").Value = '14MCDexMCB'
$ws.Range("B16").Value = 'CoinExToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D16").Value = "'0.04693"
$ws.Range("E16").Value = '15CoinExTokenCET'
$ws.Range("B17").Value = 'One'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D17").Value = "'0.0005988"
$ws.Range("E17").Value = '16OneONE'
$ws.Range("D18").Value = "'0.006279"
$ws.Range("B19").Value = 'BitKan'
$ws.Range("C19").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D19").Value = "'0.001051"
$ws.Range("E19").Value = '18BitKanKAN'
$ws.Range("B20").Value = 'HotbitToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D20").Value = "'0.003803"
$ws.Range("E20").Value = '19HotbitTokenHTB'
$ws.Range("B21").Value = 'NitroEx'
$ws.Range("C21").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D21").Value = "'0.0001504"
$ws.Range("E21").Value = '20NitroExNTX'
$ws.Range("B22").Value = 'UpBots'
$ws.Range("C22").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D22").Value = "'0.0003611"
$ws.Range("E22").Value = '21UpBotsUBXT'
$ws.Range("D23").Value = "'3.983"
$ws.Range("B24").Value = 'GateToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D24").Value = "'3.406"
$ws.Range("E24").Value = '23GateTokenGT'
$ws.Range("B25").Value = 'BTSEToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D25").Value = "'2.113"
$ws.Range("E25").Value = '24BTSETokenBTSE'
$ws.Range("E27").Value = '26ProBitTokenPROBBestin24h'
$ws.Range("D40").Value = "'0.04164"
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D41").Value = "'0.1039"
$ws.Range("E41").Value = '40BKEXTokenBKK'
$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D42").Value = "'0.003158"
$ws.Range("E42").Value = '41CEJICEJI'
$ws.Range("B43").Value = 'KickToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D43").Value = "'0.003256"
$ws.Range("E43").Value = '42KickTokenKICKWorstin24h'
$ws.Range("D44").Value = "'0.01024"
$ws.Range("D45").Value = "'0.00005651"
$ws.Range("D46").Value = "'0.00000000752"
$ws.Range("D47").Value = "'0.6821"
$ws.Range("D48").Value = "'0.02651"
$ws.Range("E48").Value = '47BOLOBOLO'
$ws.Range("D49").Value = "'0.00002106"
$ws.Range("D50").Value = "'0.01013"
